# Updated cryptos list with GitHub Actions
# Refresh price / 1h-volume-change figures pulled from coinranking.com,
# plus the FirstDigitalUSD / Aave ranking swap (rows 43 & 44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.221.69"
$ws.Range("E2").Value = "  -1.53%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.582.96"
$ws.Range("E3").Value = "  -2.21%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'563.62"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'142.66"
$ws.Range("E6").Value = "  -1.94%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.94%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.589.89"
$ws.Range("E9").Value = "  -2.31%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'6.64"
$ws.Range("E10").Value = "  -2.96%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.26%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.161"
$ws.Range("E12").Value = "  +12.20%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "'0.347"
$ws.Range("E13").Value = "  +1.59%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.039.18"
$ws.Range("E14").Value = "  -1.83%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "59.184.01"
$ws.Range("E15").Value = "  -1.58%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "'22.89"
$ws.Range("E16").Value = "  +5.42%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +0.90%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.590.94"
$ws.Range("E18").Value = "  -2.16%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'4.55"
$ws.Range("E19").Value = "  -0.71%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'337.17"
$ws.Range("E20").Value = "  -1.88%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  +0.11%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.24%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.08%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'64.24"
$ws.Range("E24").Value = "  -3.17%  "

# Row 25 - Polygon
$ws.Range("D25").Value = "'0.463"
$ws.Range("E25").Value = "  +6.12%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.18%  "

# Row 27 - Kaspa
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -2.87%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'7.32"
$ws.Range("E28").Value = "  +0.26%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0775"
$ws.Range("E29").Value = "  +0.46%  "

# Row 30 - USDe
$ws.Range("E30").Value = "  +0.08%  "

# Row 31 - Monero
$ws.Range("D31").Value = "'161.43"
$ws.Range("E31").Value = "  +3.13%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -2.21%  "

# Row 33 - Aptos
$ws.Range("D33").Value = "'6.09"
$ws.Range("E33").Value = "  -0.53%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "'18.94"
$ws.Range("E34").Value = "  -1.25%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -1.78%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -0.11%  "

# Row 37 - SuiNetwork
$ws.Range("D37").Value = "'0.875"
$ws.Range("E37").Value = "  -4.03%  "

# Row 38 - Fetch.AI
$ws.Range("D38").Value = "'0.878"
$ws.Range("E38").Value = "  -3.90%  "

# Row 39 - OKB
$ws.Range("D39").Value = "'37.47"
$ws.Range("E39").Value = "  -0.04%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -1.28%  "

# Row 41 - Bittensor
$ws.Range("D41").Value = "'295.06"
$ws.Range("E41").Value = "  -2.66%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  -0.03%  "

# Row 43 & 44 - FirstDigitalUSD / Aave swap ranking positions
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'132.46"
$ws.Range("E43").Value = "  +7.28%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.36%  "

# Row 45 - Stellar
$ws.Range("E45").Value = "  +0.07%  "

# Row 46 - Mantle
$ws.Range("D46").Value = "'0.596"
$ws.Range("E46").Value = "  -1.30%  "

# Row 47 - Hedera
$ws.Range("D47").Value = "'0.0536"
$ws.Range("E47").Value = "  -2.01%  "

# Row 48 - WhiteBITCoin
$ws.Range("D48").Value = "'10.62"
$ws.Range("E48").Value = "  -0.06%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "'19.08"
$ws.Range("E49").Value = "  -1.04%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  -1.08%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "'18.54"
$ws.Range("E51").Value = "  +0.78%  "
